# Refined metadata to be additional tab
#
# 1) Bump the per-row "time_taken" query timestamps on the existing "data"
#    sheet (F2:F8).
# 2) Add a new "metadata" worksheet (placed right after "data") that records
#    the panel query metadata (data_name/data_id/data_version/... +
#    panel_query_time/panel_get_request).

$wb = $excel.ActiveWorkbook
$data = $wb.Worksheets.Item("data")

# --- 1) refresh the time_taken column on "data" --------------------------
$data.Range("F2").Value = "2021-10-05 14:19:39.828124"
$data.Range("F3").Value = "2021-10-05 14:19:39.828133"
$data.Range("F4").Value = "2021-10-05 14:19:39.828136"
$data.Range("F5").Value = "2021-10-05 14:19:39.828139"
$data.Range("F6").Value = "2021-10-05 14:19:39.828142"
$data.Range("F7").Value = "2021-10-05 14:19:39.828145"
$data.Range("F8").Value = "2021-10-05 14:19:39.828148"

# --- 2) add the "metadata" sheet right after "data" -----------------------
$meta = $wb.Worksheets.Add($null, $data)
$meta.Name = "metadata"

# Copy the header formatting (bold / border / centered) from the "data"
# sheet's own header row so the new header row matches it exactly, then
# overwrite the copied header text with the metadata column names.
$data.Range("B1:F1").Copy()
$meta.Range("B1:F1").PasteSpecial(-4122)
$data.Range("F1").Copy()
$meta.Range("G1").PasteSpecial(-4122)

$meta.Range("B1").Value = "data_name"
$meta.Range("C1").Value = "data_id"
$meta.Range("D1").Value = "data_version"
$meta.Range("E1").Value = "data_version_created"
$meta.Range("F1").Value = "panel_query_time"
$meta.Range("G1").Value = "panel_get_request"

# Row 2: the actual metadata values, copying the index-column style from
# "data"'s A2 cell for the new sheet's A2.
$data.Range("A2").Copy()
$meta.Range("A2").PasteSpecial(-4122)
$meta.Range("A2").Value = 0

$meta.Range("B2").Value = "Common craniosynostosis syndromes"
$meta.Range("C2").Value = 507
$meta.Range("D2").Value = "'1.13"
$meta.Range("E2").Value = "2021-03-10T18:58:47.933799Z"
$meta.Range("F2").Value = "2021-10-05 14:19:39.824250"
$meta.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/507/?format=json"

# Leave the focus/selection back on "data", matching the original workbook.
$data.Activate()
